$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 3 for the columns that differ:
# D (Fecha), J (Volumen), K (Precio mínimo), L (Precio máximo),
# M (Precio promedio ponderado), O (Origen), P (Precio $/Kg)

$cols = @("D", "J", "K", "L", "M", "O", "P")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range("$col" + "2")
    $cellRow3 = $ws.Range("$col" + "3")

    $val2 = $cellRow2.Value2
    $val3 = $cellRow3.Value2

    $cellRow2.Value2 = $val3
    $cellRow3.Value2 = $val2
}
